$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the text (shared-string) columns of both data rows so the stale
# "ECs"/"MuSCs"/"Tnfsf8"/"Tnfrsf8" entries are fully released from the
# shared-strings table before we rewrite them in the desired order.
$ws.Range("A2:D3").ClearContents()

# The old row 2 (Sending cluster = ECs) is dropped; row 3 collapses up and
# becomes the new (only) data row, now re-labelled with the updated TPM
# figures.
$ws.Rows(3).Delete()

# Re-populate the label columns in the order they should appear in the
# rebuilt shared-strings table: MuSCs, Tnfsf8, Tnfrsf8, ECs.
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Tnfsf8"
$ws.Range("C2").Value = "Tnfrsf8"
$ws.Range("D2").Value = "ECs"

# Updated (new-TPM) numeric values for the remaining row.
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3050356666666666
$ws.Range("H2").Value = 0.915107
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4778466666666667
$ws.Range("N2").Value = 1.43354
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1457602765311111
$ws.Range("R2").Value = 1.31184248878
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
